# Applies the "removing 'Total' from labels" edit described in the commit:
#  - Renames a set of RACE category labels (column C) across all sheets,
#    mostly removing "Total " / " alone" wording, plus a few acronym/wording
#    tweaks.
#  - Because several of the underlying categories effectively shift meaning
#    by one slot (old "Multiracial incl. Asian" -> new "MNAW", etc.), the
#    numeric data that goes with the "detail" sheet's 4 multiracial rows per
#    region, and the "dichot" sheet's 2 multiracial rows per region, has to
#    move together with the relabeling so that each row's figures keep
#    matching what its (new) label actually means.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Label text renames (column C, every sheet)
# ---------------------------------------------------------------------------
$labelMap = @{
    "American Indian or Alaskan Native alone" = "American Indian or Alaskan Native"
    "Asian alone" = "Asian"
    "Black or African American alone" = "Black or African American"
    "Native Hawaiian and Other Pacific Islander alone" = "Native Hawaiian or Pacific Islander"
    "Some Other Race alone" = "Some Other Race"
    "White alone" = "White"
    "Multiracial incl. Asian" = "MNAW"
    "Multiracial incl. Asian, white" = "Multirace incl. Asian"
    "Multiracial incl. white" = "Multirace incl. Asian, white"
    "Multiracial not Asian or white" = "Multirace incl. white"
    "Total Multirace PSRC" = "Multirace PSRC"
    "Total Single race PSRC" = "Single race PSRC"
    "Total Multirace Harvard" = "Multirace Harvard"
    "Total Single race Harvard" = "Single race Harvard"
    "Total People of color" = "People of color"
    "Multiracial not white" = "MNW"
    "Multiple Races" = "Multirace"
}

foreach ($ws in $wb.Worksheets) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        $cur = $cell.Value()
        if ($null -ne $cur -and $labelMap.ContainsKey($cur)) {
            $cell.Value = $labelMap[$cur]
        }
    }
}

# ---------------------------------------------------------------------------
# 2) "detail" sheet: the 4-row multiracial blocks (RACE = old idx 33-36)
#    rotate their data (ARACE/HRACE counts, shares, MOEs, reliability) down
#    by one row within each block, wrapping the last row back to the first.
#    Columns: F,G (counts) I,J (shares) L,M (count moe) O,P (share moe)
#    R,S (reliability).
# ---------------------------------------------------------------------------
$detail = $wb.Worksheets.Item("detail")
$rotCols = @(6, 7, 9, 10, 12, 13, 15, 16, 18, 19)   # F G I J L M O P R S
$detailBlocks = @(
    @(11, 12, 13, 14),
    @(29, 30, 31, 32),
    @(47, 48, 49, 50),
    @(65, 66, 67, 68),
    @(83, 84, 85, 86)
)

foreach ($block in $detailBlocks) {
    $n = $block.Length
    # snapshot all the "before" values for the block first
    $snapshot = @{}
    foreach ($row in $block) {
        $rowVals = @{}
        foreach ($c in $rotCols) {
            $rowVals[$c] = $detail.Cells.Item($row, $c).Value()
        }
        $snapshot[$row] = $rowVals
    }
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $block[$i]
        $srcRow = $block[(($i - 1) + $n) % $n]
        foreach ($c in $rotCols) {
            $detail.Cells.Item($destRow, $c).Value = $snapshot[$srcRow][$c]
        }
    }
}

# ---------------------------------------------------------------------------
# 3) "dichot" sheet: the 2-row multiracial pairs (RACE = old idx 35 then 50)
#    swap their data between the two rows; the first row becomes the "MNW"
#    category and the second row becomes the "Multirace incl. white"
#    category (post-rename text).
# ---------------------------------------------------------------------------
$dichot = $wb.Worksheets.Item("dichot")
$dichotPairs = @(
    @(11, 12),
    @(27, 28),
    @(43, 44),
    @(59, 60),
    @(75, 76)
)

foreach ($pair in $dichotPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $valsA = @{}
    $valsB = @{}
    foreach ($c in $rotCols) {
        $valsA[$c] = $dichot.Cells.Item($rowA, $c).Value()
        $valsB[$c] = $dichot.Cells.Item($rowB, $c).Value()
    }

    foreach ($c in $rotCols) {
        $dichot.Cells.Item($rowA, $c).Value = $valsB[$c]
        $dichot.Cells.Item($rowB, $c).Value = $valsA[$c]
    }

    $dichot.Cells.Item($rowA, 3).Value = "MNW"
    $dichot.Cells.Item($rowB, 3).Value = "Multirace incl. white"
}
